# Update the date labels in column A from 2023-09-18..2023-09-23
# to 2023-09-25..2023-09-30 (shift each block of rows forward by one week).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of old date text -> new date text, in the order they appear top to bottom.
$oldDates = @("2023-09-18", "2023-09-19", "2023-09-20", "2023-09-21", "2023-09-22", "2023-09-23")
$newDates = @("2023-09-25", "2023-09-26", "2023-09-27", "2023-09-28", "2023-09-29", "2023-09-30")

$usedRange = $ws.UsedRange
$maxRow = $usedRange.Rows.Count

for ($r = 2; $r -le $maxRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    for ($i = 0; $i -lt $oldDates.Length; $i++) {
        if ($val -eq $oldDates[$i]) {
            $cell.Value = $newDates[$i]
            break
        }
    }
}
